# Update countries & provincias Spain
# Refresh timestamp, swap Rusia/Iran ranking and Afganistan/Islandia/Camerun
# ranking (case counts overtook one another), and update the underlying
# case numbers for the affected countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp label
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 10:22"

# Rusia overtakes Iran (row 11 <-> row 12)
$ws.Range("A11").Value = "Rusia"
$ws.Range("B11").Value = 93558
$ws.Range("C11").Value = 6411
$ws.Range("D11").Value = 8456
$ws.Range("E11").Value = 84235
$ws.Range("F11").Value = 2300
$ws.Range("G11").Value = 73
$ws.Range("H11").Value = 867

$ws.Range("A12").Value = "Iran"
$ws.Range("B12").Value = 91472
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 70933
$ws.Range("E12").Value = 14733
$ws.Range("F12").Value = 3011
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 5806

# Polonia (row 33) updated case numbers
$ws.Range("B33").Value = 12089
$ws.Range("C33").Value = 187
$ws.Range("E33").Value = 8864
$ws.Range("G33").Value = 8
$ws.Range("H33").Value = 570

# Dinamarca (row 41) updated case numbers
$ws.Range("B41").Value = 8851
$ws.Range("C41").Value = 153
$ws.Range("E41").Value = 2465

# Filipinas (row 43) updated case numbers
$ws.Range("B43").Value = 7958
$ws.Range("C43").Value = 181
$ws.Range("D43").Value = 975
$ws.Range("E43").Value = 6453
$ws.Range("G43").Value = 19
$ws.Range("H43").Value = 530

# Uzbekistan (row 68) updated case numbers
$ws.Range("D68").Value = 900
$ws.Range("E68").Value = 1016

# Afganistan overtakes Islandia and Camerun (rows 71-73 shift)
$ws.Range("A71").Value = "Afganistan"
$ws.Range("B71").Value = 1828
$ws.Range("C71").Value = 125
$ws.Range("D71").Value = 228
$ws.Range("E71").Value = 1542
$ws.Range("F71").Value = 7
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 58

$ws.Range("A72").Value = "Islandia"
$ws.Range("B72").Value = 1792
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 1624
$ws.Range("E72").Value = 158
$ws.Range("F72").Value = 1
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 10

$ws.Range("A73").Value = "Camerun"
$ws.Range("B73").Value = 1705
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 805
$ws.Range("E73").Value = 842
$ws.Range("F73").Value = 12
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 58
